$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: headers (technical codes -> human readable labels)
$ws.Range("A1").Value = "Nº hogares"
$ws.Range("B1").Value = "Tipo de hogar, código"
$ws.Range("C1").Value = "Provincia:"
$ws.Range("D1").Value = "Comarca nombre"
$ws.Range("E1").Value = "Comarca"
$ws.Range("F1").Value = "Tipo de hogar"
$ws.Range("G1").Value = "Provincia nombre"

# Row 2: measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:n-hogares"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "iaest-measure:provincia"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:comarca"
$ws.Range("F2").Value = "iaest-measure:tipo-de-hogar"
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 3: medida/dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "dim"

# Row 4: data types / URIs
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-comarca"
$ws.Range("E4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "URI-Provincia"
